$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New effect reference rows to append below the existing data (rows 1-15 already present)
$data = @(
  @(15, "Shield"),
  @(16, "Bounce Heal"),
  @(17, "Unnamed AOE Bubble"),
  @(18, "Regeneration"),
  @(19, "Chloroplast"),
  @(20, "Nature's Touch"),
  @(21, "Replenishing Winds"),
  @(22, "Blessing of the Grove Aura"),
  @(23, "Blessing of the Grove HOT"),
  @(24, "Blessing of the Grove Ref"),
  @(25, "Replenish"),
  @(26, "Living Seed"),
  @(27, "Hibernate Friend"),
  @(28, "Hibernate Foe")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = 16 + $i
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Widen column B to fit the longer new entries
$ws.Columns.Item(2).ColumnWidth = 24.7109375

# Update the view: scroll so row 7 is at the top and select the next empty row
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A30").Select() | Out-Null
